$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix simple_name values: CO2 -> co2, CH4 -> ch4
$ws.Range("A5").Value = "co2"
$ws.Range("A6").Value = "ch4"

# Fix units column: replace LaTeX-escaped per-mille code with the literal unicode character
$ws.Range("C2").Value = "‰"
$ws.Range("C3").Value = "‰"
$ws.Range("C4").Value = "‰"

# Update the selected cell to match the saved selection in the workbook
$ws.Range("C4").Select()
